$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the selection to match the used range (A2:B116) instead of full columns
$ws.Range("A2:B116").Select()

# New accuracy values for column B (row => value)
$bUpdates = @{
    3   = 0.796875
    5   = 0.71875
    6   = 0.703125
    7   = 0.65625
    8   = 0.640625
    9   = 0.609375
    10  = 0.609375
    11  = 0.578125
    12  = 0.578125
    13  = 0.609375
    15  = 0.59375
    17  = 0.609375
    19  = 0.609375
    20  = 0.609375
    21  = 0.609375
    22  = 0.609375
    23  = 0.59375
    24  = 0.59375
    25  = 0.578125
    26  = 0.578125
    27  = 0.578125
    28  = 0.578125
    29  = 0.578125
    30  = 0.578125
    31  = 0.578125
    32  = 0.578125
    33  = 0.5625
    34  = 0.5625
    35  = 0.5625
    36  = 0.5625
    37  = 0.5625
    38  = 0.5625
    39  = 0.5625
    40  = 0.5625
    41  = 0.5625
    42  = 0.5625
    43  = 0.5625
    44  = 0.5625
    45  = 0.5625
    46  = 0.5625
    47  = 0.5625
    58  = 0.578125
    59  = 0.578125
    60  = 0.578125
    61  = 0.578125
    62  = 0.578125
    63  = 0.578125
    64  = 0.578125
    65  = 0.578125
    66  = 0.578125
    91  = 0.578125
    92  = 0.578125
    93  = 0.578125
    94  = 0.578125
    95  = 0.578125
    96  = 0.578125
    97  = 0.578125
    98  = 0.578125
    99  = 0.578125
    100 = 0.578125
    101 = 0.578125
    102 = 0.578125
    103 = 0.578125
    104 = 0.5625
    105 = 0.59375
    106 = 0.671875
    107 = 0.578125
    108 = 0.53125
    109 = 0.5
    110 = 0.65625
    111 = 0.546875
    112 = 0.515625
    113 = 0.6875
    114 = 0.546875
    115 = 0.625
    116 = 0.5882352941176471
}

foreach ($row in $bUpdates.Keys) {
    $ws.Cells.Item($row, 2).Value = $bUpdates[$row]
}

# The repr() memory address embedded in column A (rows 102-116) changed between runs
$oldAddr = "0x7f79415d0e80"
$newAddr = "0x7fd76c6446d0"
for ($row = 102; $row -le 116; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $text = $cell.Value2
    $cell.Value = $text.Replace($oldAddr, $newAddr)
}
